$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added at the top of the Ciboulette block
# (rows 86-112). Insert a new row at 86, pushing the existing 86-112 block
# down to 87-113, then populate the new row 86 with the latest record.
$ws.Rows.Item(86).Insert()

$ws.Range("A86").Value = 4
$ws.Range("B86").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C86").Value = "Los Lagos"
$ws.Range("D86").Value = 44463
$ws.Range("D86").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E86").Value = 10
$ws.Range("F86").Value = 100112039
$ws.Range("G86").Value = "Ciboulette"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 240
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = 4000
$ws.Range("N86").Value = "$/docena de atados"
$ws.Range("O86").Value = "Región Metropolitana"
$ws.Range("P86").Value = 1333
$ws.Range("Q86").Value = 3
$ws.Range("R86").Value = "Hortaliza"
